$wb = $excel.ActiveWorkbook

# --- Sheet: 模型CV-95.00%Importance ---
$ws2 = $wb.Worksheets.Item("模型CV-95.00%Importance")
$ws2.Range("B2").Value = 12
$ws2.Range("B5").Value = 17
$ws2.Range("B6").Value = 8
$ws2.Range("B7").Value = 13
$ws2.Range("B8").Value = 6
$ws2.Range("B10").Value = 18
$ws2.Range("B11").Value = 14
$ws2.Range("B12").Value = 11

# --- Sheet: 特征IncreaseCV-Filter评估结果 ---
$ws3 = $wb.Worksheets.Item("特征IncreaseCV-Filter评估结果")
$ws3.Range("B2").Value = 12
$ws3.Range("B4").Value = 17
$ws3.Range("B5").Value = 8
$ws3.Range("B6").Value = 13
$ws3.Range("B7").Value = 6
$ws3.Range("B8").Value = 18
